$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40,8).Value = 0
$ws.Cells.Item(40,9).Value = 0
$ws.Cells.Item(40,11).Value = 0
$ws.Cells.Item(40,13).Value = ""
$ws.Cells.Item(113,8).Value = 10000
$ws.Cells.Item(113,9).Value = 10000
$ws.Cells.Item(113,11).Value = 10000
$ws.Cells.Item(113,13).Value = -6746
$ws.Cells.Item(135,8).Value = 696.8
$ws.Cells.Item(135,9).Value = 501.25
$ws.Cells.Item(135,11).Value = 4511.25
$ws.Cells.Item(135,13).Value = -1976.25
$ws.Cells.Item(137,8).Value = 2296.25
$ws.Cells.Item(137,9).Value = 1990
$ws.Cells.Item(137,10).Value = 2398.3333
$ws.Cells.Item(137,11).Value = 5970
$ws.Cells.Item(137,12).Value = 7194.999899999999
$ws.Cells.Item(137,13).Value = -3420
$ws.Cells.Item(137,14).Value = -12294.9999
$ws.Cells.Item(141,8).Value = 2863.074
$ws.Cells.Item(141,9).Value = 2179.4167
$ws.Cells.Item(141,10).Value = 8332.333000000001
$ws.Cells.Item(141,11).Value = 6538.250100000001
$ws.Cells.Item(141,12).Value = 24996.999
$ws.Cells.Item(141,13).Value = -1358.250100000001
$ws.Cells.Item(141,14).Value = -35356.999
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(28,8).Value = 17218.875
$ws.Cells.Item(28,10).Value = 33896.668
$ws.Cells.Item(28,12).Value = 33896.668
$ws.Cells.Item(28,14).Value = -34280.668
$ws.Cells.Item(32,8).Value = 4484.484
$ws.Cells.Item(32,9).Value = 4484.484
$ws.Cells.Item(32,11).Value = 4484.484
$ws.Cells.Item(32,13).Value = -4197.484
$ws.Cells.Item(45,8).Value = 1992.5555
$ws.Cells.Item(45,9).Value = 1739.3334
$ws.Cells.Item(45,10).Value = 2499
$ws.Cells.Item(45,11).Value = 1739.3334
$ws.Cells.Item(45,12).Value = 2499
$ws.Cells.Item(45,13).Value = -1362.3334
$ws.Cells.Item(45,14).Value = -3253
$ws.Cells.Item(61,8).Value = 3470.1667
$ws.Cells.Item(61,9).Value = 2956
$ws.Cells.Item(61,11).Value = 2956
$ws.Cells.Item(61,13).Value = -2744
$ws.Cells.Item(74,8).Value = 5634
$ws.Cells.Item(74,9).Value = 4451
$ws.Cells.Item(74,10).Value = 8000
$ws.Cells.Item(74,11).Value = 4451
$ws.Cells.Item(74,12).Value = 8000
$ws.Cells.Item(74,13).Value = -3577
$ws.Cells.Item(74,14).Value = -9748
$ws.Cells.Item(77,8).Value = 5634
$ws.Cells.Item(77,9).Value = 4451
$ws.Cells.Item(77,10).Value = 8000
$ws.Cells.Item(77,11).Value = 22255
$ws.Cells.Item(77,12).Value = 40000
$ws.Cells.Item(77,13).Value = -17887
$ws.Cells.Item(77,14).Value = -48736
$ws.Cells.Item(99,8).Value = 17218.875
$ws.Cells.Item(99,10).Value = 33896.668
$ws.Cells.Item(99,12).Value = 33896.668
$ws.Cells.Item(99,14).Value = -39886.668
$ws.Cells.Item(122,8).Value = 0
$ws.Cells.Item(122,9).Value = 0
$ws.Cells.Item(122,10).Value = 0
$ws.Cells.Item(122,11).Value = 0
$ws.Cells.Item(122,12).Value = 0
$ws.Cells.Item(122,13).Value = ""
$ws.Cells.Item(122,14).Value = ""
$ws.Cells.Item(136,8).Value = 3470.1667
$ws.Cells.Item(136,9).Value = 2956
$ws.Cells.Item(136,11).Value = 8868
$ws.Cells.Item(136,13).Value = -6318
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86,8).Value = 952.25
$ws.Cells.Item(86,9).Value = 1539.8334
$ws.Cells.Item(86,10).Value = 599.7
$ws.Cells.Item(86,11).Value = 1539.8334
$ws.Cells.Item(86,12).Value = 599.7
$ws.Cells.Item(86,13).Value = -416.8334
$ws.Cells.Item(86,14).Value = -2845.7
$ws.Cells.Item(89,8).Value = 952.25
$ws.Cells.Item(89,9).Value = 1539.8334
$ws.Cells.Item(89,10).Value = 599.7
$ws.Cells.Item(89,11).Value = 7699.166999999999
$ws.Cells.Item(89,12).Value = 2998.5
$ws.Cells.Item(89,13).Value = -2083.166999999999
$ws.Cells.Item(89,14).Value = -14230.5
$ws.Cells.Item(102,8).Value = 45000
$ws.Cells.Item(102,9).Value = 45000
$ws.Cells.Item(102,10).Value = 0
$ws.Cells.Item(102,11).Value = 45000
$ws.Cells.Item(102,12).Value = 0
$ws.Cells.Item(102,13).Value = -41755
$ws.Cells.Item(102,14).Value = ""
$ws.Cells.Item(134,8).Value = 3902.1714
$ws.Cells.Item(134,9).Value = 4065.6
$ws.Cells.Item(134,10).Value = 2921.6
$ws.Cells.Item(134,11).Value = 12196.8
$ws.Cells.Item(134,12).Value = 8764.799999999999
$ws.Cells.Item(134,13).Value = -9661.799999999999
$ws.Cells.Item(134,14).Value = -13834.8
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(13,8).Value = 5001.5
$ws.Cells.Item(13,10).Value = 9999
$ws.Cells.Item(13,12).Value = 9999
$ws.Cells.Item(13,14).Value = -10277
$ws.Cells.Item(16,8).Value = 545.7778
$ws.Cells.Item(16,9).Value = 627.1667
$ws.Cells.Item(16,11).Value = 627.1667
$ws.Cells.Item(16,13).Value = -340.1667
$ws.Cells.Item(62,8).Value = 4833.3335
$ws.Cells.Item(62,9).Value = 4250
$ws.Cells.Item(62,11).Value = 4250
$ws.Cells.Item(62,13).Value = -3626
$ws.Cells.Item(65,8).Value = 4833.3335
$ws.Cells.Item(65,9).Value = 4250
$ws.Cells.Item(65,11).Value = 21250
$ws.Cells.Item(65,13).Value = -18130
$ws.Cells.Item(113,8).Value = 545.7778
$ws.Cells.Item(113,9).Value = 627.1667
$ws.Cells.Item(113,11).Value = 627.1667
$ws.Cells.Item(113,13).Value = 1542.8333
$ws.Cells.Item(132,8).Value = 2185.4
$ws.Cells.Item(132,9).Value = 980
$ws.Cells.Item(132,10).Value = 3993.5
$ws.Cells.Item(132,11).Value = 2940
$ws.Cells.Item(132,12).Value = 11980.5
$ws.Cells.Item(132,13).Value = -410
$ws.Cells.Item(132,14).Value = -17040.5
$ws.Cells.Item(134,8).Value = 1958.16
$ws.Cells.Item(134,9).Value = 1830.8334
$ws.Cells.Item(134,11).Value = 5492.5002
$ws.Cells.Item(134,13).Value = -2957.5002
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23,8).Value = 341.85715
$ws.Cells.Item(23,10).Value = 361.4
$ws.Cells.Item(23,12).Value = 1084.2
$ws.Cells.Item(23,14).Value = -1554.2
$ws.Cells.Item(33,8).Value = 119.833336
$ws.Cells.Item(33,9).Value = 73
$ws.Cells.Item(33,11).Value = 438
$ws.Cells.Item(33,13).Value = -155
$ws.Cells.Item(75,8).Value = 4067.2222
$ws.Cells.Item(75,10).Value = 4184.8335
$ws.Cells.Item(75,12).Value = 12554.5005
$ws.Cells.Item(75,14).Value = -14550.5005
$ws.Cells.Item(78,8).Value = 4067.2222
$ws.Cells.Item(78,10).Value = 4184.8335
$ws.Cells.Item(78,12).Value = 37663.5015
$ws.Cells.Item(78,14).Value = -47647.5015
$ws.Cells.Item(80,8).Value = 6436.75
$ws.Cells.Item(80,9).Value = 1797.5
$ws.Cells.Item(80,10).Value = 7983.1665
$ws.Cells.Item(80,11).Value = 5392.5
$ws.Cells.Item(80,12).Value = 23949.4995
$ws.Cells.Item(80,13).Value = -4456.5
$ws.Cells.Item(80,14).Value = -25821.4995
$ws.Cells.Item(83,8).Value = 6436.75
$ws.Cells.Item(83,9).Value = 1797.5
$ws.Cells.Item(83,10).Value = 7983.1665
$ws.Cells.Item(83,11).Value = 16177.5
$ws.Cells.Item(83,12).Value = 71848.4985
$ws.Cells.Item(83,13).Value = -11497.5
$ws.Cells.Item(83,14).Value = -81208.4985
$ws.Cells.Item(104,8).Value = 2016.5714
$ws.Cells.Item(122,8).Value = 687.8421
$ws.Cells.Item(122,9).Value = 709.625
$ws.Cells.Item(122,10).Value = 672
$ws.Cells.Item(122,11).Value = 6386.625
$ws.Cells.Item(122,12).Value = 6048
$ws.Cells.Item(122,13).Value = -3936.625
$ws.Cells.Item(122,14).Value = -10948
$ws.Cells.Item(140,8).Value = 1669.7693
$ws.Cells.Item(140,9).Value = 1350.5834
$ws.Cells.Item(140,11).Value = 4051.7502
$ws.Cells.Item(140,13).Value = 1128.2498
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97,8).Value = 721.8889
$ws.Cells.Item(97,9).Value = 485.2857
$ws.Cells.Item(97,10).Value = 1550
$ws.Cells.Item(97,11).Value = 485.2857
$ws.Cells.Item(97,12).Value = 1550
$ws.Cells.Item(97,13).Value = 10.71429999999998
$ws.Cells.Item(97,14).Value = -2542
$ws.Cells.Item(122,8).Value = 11366902
$ws.Cells.Item(122,9).Value = 15626765
$ws.Cells.Item(122,11).Value = 46880295
$ws.Cells.Item(122,13).Value = -46877845
$ws.Cells.Item(126,8).Value = 1998.5
$ws.Cells.Item(126,9).Value = 1000
$ws.Cells.Item(126,10).Value = 2997
$ws.Cells.Item(126,11).Value = 3000
$ws.Cells.Item(126,12).Value = 8991
$ws.Cells.Item(126,13).Value = -530
$ws.Cells.Item(126,14).Value = -13931
$ws.Cells.Item(132,8).Value = 2262.923
$ws.Cells.Item(132,10).Value = 1000
$ws.Cells.Item(132,12).Value = 3000
$ws.Cells.Item(132,14).Value = -8060
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122,8).Value = 3753.5
$ws.Cells.Item(122,9).Value = 3603.2
$ws.Cells.Item(122,11).Value = 10809.6
$ws.Cells.Item(122,13).Value = -8359.599999999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(9,8).Value = 106
$ws.Cells.Item(9,9).Value = 106
$ws.Cells.Item(9,10).Value = 0
$ws.Cells.Item(9,11).Value = 106
$ws.Cells.Item(9,12).Value = 0
$ws.Cells.Item(9,13).Value = 34
$ws.Cells.Item(9,14).Value = ""
$ws.Cells.Item(14,8).Value = 3697.1428
$ws.Cells.Item(14,9).Value = 15000
$ws.Cells.Item(14,10).Value = 1813.3334
$ws.Cells.Item(14,11).Value = 15000
$ws.Cells.Item(14,12).Value = 1813.3334
$ws.Cells.Item(14,13).Value = -14832
$ws.Cells.Item(14,14).Value = -2149.3334
$ws.Cells.Item(20,8).Value = 3239.5
$ws.Cells.Item(20,9).Value = 1468
$ws.Cells.Item(20,11).Value = 1468
$ws.Cells.Item(20,13).Value = -1228
$ws.Cells.Item(54,8).Value = 26000
$ws.Cells.Item(54,10).Value = 29500
$ws.Cells.Item(54,12).Value = 29500
$ws.Cells.Item(54,14).Value = -30540
$ws.Cells.Item(136,8).Value = 5191.533
$ws.Cells.Item(136,9).Value = 4836.385
$ws.Cells.Item(136,11).Value = 14509.155
$ws.Cells.Item(136,13).Value = -11959.155
